# Insert a brand-new record as row 303 in the single sheet, pushing the
# existing rows 303-401 down to 304-402 (dimension grows from A1:T401 to
# A1:T402). This matches a weekly refresh of the "Fruta, Feria Lagunitas
# de Puerto Montt - Plátano" price series: one new daily observation is
# prepended to the historical table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 303 (and everything below it) down by one row.
$ws.Rows.Item(303).EntireRow.Insert()

# Populate the newly freed row 303 with the new observation.
$ws.Cells.Item(303, 1).Value  = 4
$ws.Cells.Item(303, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(303, 3).Value  = "Los Lagos"
$ws.Cells.Item(303, 4).Value  = 44627
$ws.Cells.Item(303, 5).Value  = 10
$ws.Cells.Item(303, 6).Value  = "Fruta"
$ws.Cells.Item(303, 7).Value  = 100108
$ws.Cells.Item(303, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(303, 9).Value  = 100108006
$ws.Cells.Item(303, 10).Value = "Plátano"
$ws.Cells.Item(303, 11).Value = "Sin especificar"
$ws.Cells.Item(303, 12).Value = "Primera Pintón"
$ws.Cells.Item(303, 13).Value = 600
$ws.Cells.Item(303, 14).Value = 19000
$ws.Cells.Item(303, 15).Value = 19000
$ws.Cells.Item(303, 16).Value = 19000
$ws.Cells.Item(303, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(303, 18).Value = "Ecuador"
$ws.Cells.Item(303, 19).Value = 950
$ws.Cells.Item(303, 20).Value = 20
